$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country label reassignments (caused by re-ranking rows after new case counts) ---
# Shared-string table order changed so that these rows now display different country
# names even though the row position (rank by "Casos totales") stays the same.
$ws.Range("A96").Value = "Mauritania"
$ws.Range("A97").Value = "Somalia"
$ws.Range("A98").Value = "Republica de Africa Central"
$ws.Range("A169").Value = "Angola"
$ws.Range("A170").Value = "Gibraltar"
$ws.Range("A171").Value = "Guadalupe"
$ws.Range("A202").Value = "Fiyi"
$ws.Range("A203").Value = "Dominica"
$ws.Range("A208").Value = "Santa Sede"
$ws.Range("A209").Value = "Islas Turcas y Caicos"
$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("A214").Value = "Islas Virgenes Britanicas"

# --- Updated statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes) ---
$ws.Range("B4").Value = 2318630
$ws.Range("C4").Value = 21440
$ws.Range("D4").Value = 962880
$ws.Range("E4").Value = 1234000
$ws.Range("G4").Value = 343
$ws.Range("H4").Value = 121750
$ws.Range("B7").Value = 411478
$ws.Range("C7").Value = 15666
$ws.Range("D7").Value = 228158
$ws.Range("E7").Value = 169974
$ws.Range("G7").Value = 376
$ws.Range("H7").Value = 13346
$ws.Range("D14").Value = 174700
$ws.Range("E14").Value = 7305
$ws.Range("B18").Value = 160093
$ws.Range("C18").Value = 641
$ws.Range("E18").Value = 56343
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = 29633
$ws.Range("B22").Value = 92681
$ws.Range("C22").Value = 4966
$ws.Range("D22").Value = 50326
$ws.Range("E22").Value = 40478
$ws.Range("G22").Value = 46
$ws.Range("H22").Value = 1877
$ws.Range("B29").Value = 53758
$ws.Range("C29").Value = 1547
$ws.Range("D29").Value = 14327
$ws.Range("E29").Value = 37325
$ws.Range("G29").Value = 89
$ws.Range("H29").Value = 2106
$ws.Range("B58").Value = 13717
$ws.Range("C58").Value = 514
$ws.Range("D58").Value = 10074
$ws.Range("E58").Value = 3558
$ws.Range("G58").Value = 15
$ws.Range("H58").Value = 85
$ws.Range("B76").Value = 6153
$ws.Range("C76").Value = 207
$ws.Range("E76").Value = 1844
$ws.Range("B96").Value = 2813
$ws.Range("C96").Value = 192
$ws.Range("D96").Value = 696
$ws.Range("E96").Value = 2009
$ws.Range("G96").Value = 6
$ws.Range("H96").Value = 108
$ws.Range("B97").Value = 2755
$ws.Range("C97").Value = 36
$ws.Range("D97").Value = 751
$ws.Range("E97").Value = 1916
$ws.Range("H97").Value = 88
$ws.Range("B98").Value = 2686
$ws.Range("C98").Value = 81
$ws.Range("D98").Value = 420
$ws.Range("E98").Value = 2247
$ws.Range("H98").Value = 19
$ws.Range("B104").Value = 2127
$ws.Range("C104").Value = 69
$ws.Range("D104").Value = 1014
$ws.Range("E104").Value = 1101
$ws.Range("B121").Value = 1362
$ws.Range("C121").Value = 26
$ws.Range("D121").Value = 791
$ws.Range("E121").Value = 558
$ws.Range("B137").Value = 783
$ws.Range("C137").Value = 108
$ws.Range("E137").Value = 343
$ws.Range("B147").Value = 627
$ws.Range("C147").Value = 4
$ws.Range("D147").Value = 285
$ws.Range("E147").Value = 337
$ws.Range("G147").Value = 1
$ws.Range("H147").Value = 5
$ws.Range("C169").Value = 4
$ws.Range("D169").Value = 66
$ws.Range("E169").Value = 101
$ws.Range("G169").Value = 1
$ws.Range("H169").Value = 9
$ws.Range("B170").Value = 176
$ws.Range("D170").Value = 176
$ws.Range("E170").Value = 0
$ws.Range("H170").Value = 0
$ws.Range("B171").Value = 174
$ws.Range("D171").Value = 157
$ws.Range("E171").Value = 3
$ws.Range("H171").Value = 14
$ws.Range("D208").Value = 12
$ws.Range("H208").Value = 0
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 1
$ws.Range("D213").Value = 8
$ws.Range("H213").Value = 0
$ws.Range("D214").Value = 7
$ws.Range("H214").Value = 1

# --- Footer timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 20 de Junio de 2020 a las 21:15"
